# GridMaze maze_poke_1_0 BOM update:
# Poke IR beam current reduced (R1: 82R -> 90R), updating part number and
# Farnell order code accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Row 5 is R1 (the resistor that sets the IR beam current)
$ws.Range("G5").Value = "MCWR08X90R9FTL"
$ws.Range("B5").Value = "90R"
$ws.Range("H5").Value = 2695100

# Leave the selection where the author last left it when saving
$ws.Range("E30").Select() | Out-Null
